$wb = $excel.ActiveWorkbook

# This script applies an updated market-price data refresh (currentAveragePrice,
# LevePriceNQ/HQ, LeveProfitNQ/HQ columns) to specific rows across the eight
# Leve-profession worksheets, as produced by the scheduled pricing data runner.

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 183955.83
$ws.Range("J13").Value = 232.5
$ws.Range("L13").Value = 232.5
$ws.Range("N13").Value = -570.5
# Row 40
$ws.Range("H40").Value = 3556.7856
$ws.Range("I40").Value = 2970.7144
$ws.Range("J40").Value = 4142.857
$ws.Range("K40").Value = 2970.7144
$ws.Range("L40").Value = 4142.857
$ws.Range("M40").Value = -2795.7144
$ws.Range("N40").Value = -4492.857
# Row 55
$ws.Range("H55").Value = 382.58823
$ws.Range("I55").Value = 187
$ws.Range("K55").Value = 187
$ws.Range("M55").Value = 27
# Row 98
$ws.Range("H98").Value = 855.5714
$ws.Range("I98").Value = 855.5714
$ws.Range("K98").Value = 855.5714
$ws.Range("M98").Value = 642.4286
# Row 121
$ws.Range("H121").Value = 3212.25
$ws.Range("J121").Value = 3212.25
$ws.Range("L121").Value = 9636.75
$ws.Range("N121").Value = -13130.75
# Row 122
$ws.Range("H122").Value = 855.5714
$ws.Range("I122").Value = 855.5714
$ws.Range("K122").Value = 2566.7142
$ws.Range("M122").Value = -116.7142000000003
# Row 138
$ws.Range("H138").Value = 2556.926
$ws.Range("I138").Value = 1763.1818
$ws.Range("K138").Value = 5289.5454
$ws.Range("M138").Value = -149.5454

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 31385.057
$ws.Range("I32").Value = 18249.77
$ws.Range("J32").Value = 120413.11
$ws.Range("K32").Value = 18249.77
$ws.Range("L32").Value = 120413.11
$ws.Range("M32").Value = -17962.77
$ws.Range("N32").Value = -120987.11
# Row 41
$ws.Range("H41").Value = 2868.3076
$ws.Range("I41").Value = 2228.8
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 2228.8
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -1814.8
$ws.Range("N41").Value = -5828
# Row 45
$ws.Range("H45").Value = 483832.34
$ws.Range("I45").Value = 1013018.2
$ws.Range("J45").Value = 2754.2727
$ws.Range("K45").Value = 1013018.2
$ws.Range("L45").Value = 2754.2727
$ws.Range("M45").Value = -1012641.2
$ws.Range("N45").Value = -3508.2727
# Row 74
$ws.Range("H74").Value = 1950.2632
$ws.Range("I74").Value = 1816
$ws.Range("J74").Value = 2666.3333
$ws.Range("K74").Value = 1816
$ws.Range("L74").Value = 2666.3333
$ws.Range("M74").Value = -942
$ws.Range("N74").Value = -4414.3333
# Row 77
$ws.Range("H77").Value = 1950.2632
$ws.Range("I77").Value = 1816
$ws.Range("J77").Value = 2666.3333
$ws.Range("K77").Value = 9080
$ws.Range("L77").Value = 13331.6665
$ws.Range("M77").Value = -4712
$ws.Range("N77").Value = -22067.6665
# Row 101
$ws.Range("H101").Value = 40602
$ws.Range("J101").Value = 40602
$ws.Range("L101").Value = 40602
$ws.Range("N101").Value = -47092
# Row 122
$ws.Range("H122").Value = 1927.174
$ws.Range("I122").Value = 1822.4375
$ws.Range("K122").Value = 5467.3125
$ws.Range("M122").Value = -3017.3125
# Row 132
$ws.Range("H132").Value = 1685.2858
$ws.Range("J132").Value = 2645
$ws.Range("L132").Value = 7935
$ws.Range("N132").Value = -12995

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6764.125
$ws.Range("I20").Value = 6024.85
$ws.Range("K20").Value = 6024.85
$ws.Range("M20").Value = -5777.85
# Row 22
$ws.Range("H22").Value = 278
$ws.Range("I22").Value = 278
$ws.Range("K22").Value = 278
$ws.Range("M22").Value = -105
# Row 86
$ws.Range("H86").Value = 1976.04
$ws.Range("I86").Value = 1392.8572
$ws.Range("J86").Value = 2718.2727
$ws.Range("K86").Value = 1392.8572
$ws.Range("L86").Value = 2718.2727
$ws.Range("M86").Value = -269.8571999999999
$ws.Range("N86").Value = -4964.2727
# Row 89
$ws.Range("H89").Value = 1976.04
$ws.Range("I89").Value = 1392.8572
$ws.Range("J89").Value = 2718.2727
$ws.Range("K89").Value = 6964.286
$ws.Range("L89").Value = 13591.3635
$ws.Range("M89").Value = -1348.286
$ws.Range("N89").Value = -24823.3635

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1393.4348
$ws.Range("I16").Value = 1214
$ws.Range("K16").Value = 1214
$ws.Range("M16").Value = -927
# Row 31
$ws.Range("H31").Value = 1899.8125
$ws.Range("I31").Value = 1269.5555
$ws.Range("J31").Value = 3790.5833
$ws.Range("K31").Value = 1269.5555
$ws.Range("L31").Value = 3790.5833
$ws.Range("M31").Value = -974.5554999999999
$ws.Range("N31").Value = -4380.5833
# Row 34
$ws.Range("H34").Value = 1899.8125
$ws.Range("I34").Value = 1269.5555
$ws.Range("J34").Value = 3790.5833
$ws.Range("K34").Value = 1269.5555
$ws.Range("L34").Value = 3790.5833
$ws.Range("M34").Value = -1067.5555
$ws.Range("N34").Value = -4194.5833
# Row 41
$ws.Range("H41").Value = 13997.091
$ws.Range("J41").Value = 13997.091
$ws.Range("L41").Value = 13997.091
$ws.Range("N41").Value = -14853.091
# Row 51
$ws.Range("H51").Value = 14249.0625
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 14999
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 14999
$ws.Range("M51").Value = -2264
$ws.Range("N51").Value = -16471
# Row 60
$ws.Range("H60").Value = 9880.3125
$ws.Range("J60").Value = 10768.615
$ws.Range("L60").Value = 10768.615
$ws.Range("N60").Value = -11790.615
# Row 61
$ws.Range("H61").Value = 14249.0625
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 14999
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 14999
$ws.Range("M61").Value = -2652
$ws.Range("N61").Value = -15695
# Row 68
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498
# Row 71
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488
# Row 107
$ws.Range("H107").Value = 1127.5
$ws.Range("J107").Value = 824.8333
$ws.Range("L107").Value = 824.8333
$ws.Range("N107").Value = -4664.8333
# Row 113
$ws.Range("H113").Value = 1393.4348
$ws.Range("I113").Value = 1214
$ws.Range("K113").Value = 1214
$ws.Range("M113").Value = 956
# Row 122
$ws.Range("H122").Value = 7778.1113
$ws.Range("I122").Value = 7102.2
$ws.Range("K122").Value = 21306.6
$ws.Range("M122").Value = -18856.6
# Row 132
$ws.Range("H132").Value = 2656.353
$ws.Range("I132").Value = 2343.9333
$ws.Range("K132").Value = 7031.7999
$ws.Range("M132").Value = -4501.7999

$ws = $wb.Worksheets.Item("CUL")
# Row 61
$ws.Range("H61").Value = 211.53334
$ws.Range("I61").Value = 122.5
$ws.Range("J61").Value = 313.2857
$ws.Range("K61").Value = 367.5
$ws.Range("L61").Value = 939.8571000000001
$ws.Range("M61").Value = -152.5
$ws.Range("N61").Value = -1369.8571
# Row 131
$ws.Range("H131").Value = 20223.428
$ws.Range("J131").Value = 24594.77
$ws.Range("L131").Value = 73784.31
$ws.Range("N131").Value = -83864.31

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 21508
$ws.Range("I40").Value = 13016
$ws.Range("K40").Value = 13016
$ws.Range("M40").Value = -12865
# Row 80
$ws.Range("H80").Value = 6498.3335
$ws.Range("I80").Value = 3499.5
$ws.Range("J80").Value = 7997.75
$ws.Range("K80").Value = 3499.5
$ws.Range("L80").Value = 7997.75
$ws.Range("M80").Value = -2501.5
$ws.Range("N80").Value = -9993.75
# Row 83
$ws.Range("H83").Value = 6498.3335
$ws.Range("I83").Value = 3499.5
$ws.Range("J83").Value = 7997.75
$ws.Range("K83").Value = 17497.5
$ws.Range("L83").Value = 39988.75
$ws.Range("M83").Value = -12505.5
$ws.Range("N83").Value = -49972.75
# Row 86
$ws.Range("H86").Value = 19724.666
$ws.Range("J86").Value = 19724.666
$ws.Range("L86").Value = 19724.666
$ws.Range("N86").Value = -22096.666
# Row 89
$ws.Range("H89").Value = 19724.666
$ws.Range("J89").Value = 19724.666
$ws.Range("L89").Value = 59173.99800000001
$ws.Range("N89").Value = -71029.99800000001
# Row 102
$ws.Range("H102").Value = 2295.6
$ws.Range("I102").Value = 1963.0769
$ws.Range("K102").Value = 1963.0769
$ws.Range("M102").Value = -341.0769
# Row 122
$ws.Range("H122").Value = 2568.037
$ws.Range("I122").Value = 2766.318
$ws.Range("J122").Value = 1695.6
$ws.Range("K122").Value = 8298.954000000002
$ws.Range("L122").Value = 5086.799999999999
$ws.Range("M122").Value = -5848.954000000002
$ws.Range("N122").Value = -9986.799999999999
# Row 126
$ws.Range("H126").Value = 2866.913
$ws.Range("I126").Value = 2835.2144
$ws.Range("K126").Value = 8505.643199999999
$ws.Range("M126").Value = -6035.643199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 566
$ws.Range("I16").Value = 566
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 566
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -396
$ws.Range("N16").ClearContents()
# Row 122
$ws.Range("H122").Value = 18027.428
$ws.Range("I122").Value = 23378.4
$ws.Range("K122").Value = 70135.20000000001
$ws.Range("M122").Value = -67685.20000000001
# Row 132
$ws.Range("H132").Value = 321996.1
$ws.Range("I132").Value = 470756.4
$ws.Range("K132").Value = 1412269.2
$ws.Range("M132").Value = -1409739.2

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 107
$ws.Range("H107").Value = 100001380
$ws.Range("I107").Value = 1725
$ws.Range("J107").Value = 500000000
$ws.Range("K107").Value = 5175
$ws.Range("L107").Value = 1500000000
$ws.Range("M107").Value = -3255
$ws.Range("N107").Value = -1500003840
# Row 122
$ws.Range("H122").Value = 7273.6
$ws.Range("I122").Value = 7304
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 21912
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -19462
$ws.Range("N122").Value = -25900
# Row 132
$ws.Range("H132").Value = 25741.117
$ws.Range("I132").Value = 26080.857
$ws.Range("K132").Value = 78242.571
$ws.Range("M132").Value = -75712.571
